$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row: rename/rearrange header labels ---
# Columns D, E, F used to be ClientID / ProjectTypeID / EmployeeCode (all unused/empty below),
# now they are ClientId / ProjectTypeId / EmployeeId and are populated with real foreign-key ids.
# Column G keeps its header "AccountingName".
$ws.Range("D1").Value = "ClientId"
$ws.Range("E1").Value = "ProjectTypeId"
$ws.Range("F1").Value = "EmployeeId"
$ws.Range("G1").Value = "AccountingName"

# --- Data rows: fill in the ClientId / ProjectTypeId / EmployeeId columns ---
$ws.Range("D2").Value = 6
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1

$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 2

$ws.Range("D4").Value = 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 3

$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 4

$ws.Range("D6").Value = 2
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 5

$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 5
$ws.Range("F7").Value = 6

$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 10
$ws.Range("F8").Value = 7

$ws.Range("D9").Value = 4
$ws.Range("E9").Value = 11
$ws.Range("F9").Value = 8

$ws.Range("D10").Value = 3
$ws.Range("E10").Value = 5
$ws.Range("F10").Value = 9

$ws.Range("D11").Value = 2
$ws.Range("E11").Value = 11
$ws.Range("F11").Value = 10

$ws.Range("D12").Value = 3
$ws.Range("E12").Value = 11
$ws.Range("F12").Value = 11

# --- Update the active selection to match the saved view state ---
$ws.Range("F20").Select()
